# Refresh cached market-board figures (currentAveragePrice*, LevePrice*,
# LeveProfit*) for a handful of leve rows across several crafting-job
# sheets, per the scheduled market-data sync.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 2653.875
$ws.Range("I34").Value = 2653.875
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 2653.875
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -2450.875
$ws.Range("N34").ClearContents()
$ws.Range("H36").Value = 2653.875
$ws.Range("I36").Value = 2653.875
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 2653.875
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -1938.875
$ws.Range("N36").ClearContents()
$ws.Range("H76").Value = 6176939.5
$ws.Range("I76").Value = 6176939.5
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 6176939.5
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -6176624.5
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 6176939.5
$ws.Range("I79").Value = 6176939.5
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 6176939.5
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -6175847.5
$ws.Range("N79").ClearContents()
$ws.Range("H100").Value = 12822508
$ws.Range("I100").Value = 18519900
$ws.Range("J100").Value = 3373.25
$ws.Range("K100").Value = 18519900
$ws.Range("L100").Value = 3373.25
$ws.Range("M100").Value = -18519359
$ws.Range("N100").Value = -4455.25
$ws.Range("H141").Value = 3027.0386
$ws.Range("I141").Value = 2936.2727
$ws.Range("J141").Value = 3526.25
$ws.Range("K141").Value = 8808.8181
$ws.Range("L141").Value = 10578.75
$ws.Range("M141").Value = -3628.8181
$ws.Range("N141").Value = -20938.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 6006.2856
$ws.Range("I63").Value = 4009.5
$ws.Range("J63").Value = 8668.666999999999
$ws.Range("K63").Value = 4009.5
$ws.Range("L63").Value = 8668.666999999999
$ws.Range("M63").Value = -3323.5
$ws.Range("N63").Value = -10040.667
$ws.Range("H66").Value = 6006.2856
$ws.Range("I66").Value = 4009.5
$ws.Range("J66").Value = 8668.666999999999
$ws.Range("K66").Value = 20047.5
$ws.Range("L66").Value = 43343.335
$ws.Range("M66").Value = -16615.5
$ws.Range("N66").Value = -50207.335
$ws.Range("H74").Value = 1287.1666
$ws.Range("I74").Value = 954.65216
$ws.Range("J74").Value = 2379.7144
$ws.Range("K74").Value = 954.65216
$ws.Range("L74").Value = 2379.7144
$ws.Range("M74").Value = -80.65215999999998
$ws.Range("N74").Value = -4127.7144
$ws.Range("H77").Value = 1287.1666
$ws.Range("I77").Value = 954.65216
$ws.Range("J77").Value = 2379.7144
$ws.Range("K77").Value = 4773.2608
$ws.Range("L77").Value = 11898.572
$ws.Range("M77").Value = -405.2608
$ws.Range("N77").Value = -20634.572
$ws.Range("H88").Value = 4047.8235
$ws.Range("I88").Value = 2181.2
$ws.Range("K88").Value = 2181.2
$ws.Range("M88").Value = -1775.2
$ws.Range("H91").Value = 4047.8235
$ws.Range("I91").Value = 2181.2
$ws.Range("K91").Value = 2181.2
$ws.Range("M91").Value = -777.1999999999998
$ws.Range("H97").Value = 3200.111
$ws.Range("I97").Value = 2298.3333
$ws.Range("J97").Value = 5003.6665
$ws.Range("K97").Value = 2298.3333
$ws.Range("L97").Value = 5003.6665
$ws.Range("M97").Value = -1802.3333
$ws.Range("N97").Value = -5995.6665

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 3333383.2
$ws.Range("I7").Value = 100
$ws.Range("J7").Value = 5000025
$ws.Range("K7").Value = 100
$ws.Range("L7").Value = 5000025
$ws.Range("M7").Value = 13
$ws.Range("N7").Value = -5000251
$ws.Range("H30").Value = 3110
$ws.Range("I30").Value = 600
$ws.Range("J30").Value = 3612
$ws.Range("K30").Value = 600
$ws.Range("L30").Value = 3612
$ws.Range("M30").Value = -475
$ws.Range("N30").Value = -3862
$ws.Range("H86").Value = 35719110
$ws.Range("I86").Value = 52634156
$ws.Range("J86").Value = 9573.777
$ws.Range("K86").Value = 52634156
$ws.Range("L86").Value = 9573.777
$ws.Range("M86").Value = -52633033
$ws.Range("N86").Value = -11819.777
$ws.Range("H89").Value = 35719110
$ws.Range("I89").Value = 52634156
$ws.Range("J89").Value = 9573.777
$ws.Range("K89").Value = 263170780
$ws.Range("L89").Value = 47868.885
$ws.Range("M89").Value = -263165164
$ws.Range("N89").Value = -59100.885
$ws.Range("H94").Value = 1037.92
$ws.Range("I94").Value = 1068.7778
$ws.Range("J94").Value = 958.5714
$ws.Range("K94").Value = 1068.7778
$ws.Range("L94").Value = 958.5714
$ws.Range("M94").Value = -617.7778000000001
$ws.Range("N94").Value = -1860.5714
$ws.Range("H105").Value = 1605.4546
$ws.Range("I105").Value = 1538.421
$ws.Range("J105").Value = 1696.4286
$ws.Range("K105").Value = 1538.421
$ws.Range("L105").Value = 1696.4286
$ws.Range("M105").Value = 208.579
$ws.Range("N105").Value = -5190.4286
$ws.Range("H134").Value = 5637.387
$ws.Range("I134").Value = 850.3333
$ws.Range("K134").Value = 2550.9999
$ws.Range("M134").Value = -15.9998999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7144666.5
$ws.Range("I31").Value = 9616797
$ws.Range("J31").Value = 2957.5557
$ws.Range("K31").Value = 9616797
$ws.Range("L31").Value = 2957.5557
$ws.Range("M31").Value = -9616502
$ws.Range("N31").Value = -3547.5557
$ws.Range("H34").Value = 7144666.5
$ws.Range("I34").Value = 9616797
$ws.Range("J34").Value = 2957.5557
$ws.Range("K34").Value = 9616797
$ws.Range("L34").Value = 2957.5557
$ws.Range("M34").Value = -9616595
$ws.Range("N34").Value = -3361.5557
$ws.Range("H134").Value = 20607.295
$ws.Range("I134").Value = 21399.428
$ws.Range("K134").Value = 64198.284
$ws.Range("M134").Value = -61663.284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 1125.25
$ws.Range("I126").Value = 500.33334
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 1501.00002
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = 3438.99998
$ws.Range("N126").Value = -18880

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 19404200
$ws.Range("I70").Value = 40185972
$ws.Range("J70").Value = 7880.467
$ws.Range("K70").Value = 40185972
$ws.Range("L70").Value = 7880.467
$ws.Range("M70").Value = -40185702
$ws.Range("N70").Value = -8420.467000000001
$ws.Range("H73").Value = 19404200
$ws.Range("I73").Value = 40185972
$ws.Range("J73").Value = 7880.467
$ws.Range("K73").Value = 40185972
$ws.Range("L73").Value = 7880.467
$ws.Range("M73").Value = -40185036
$ws.Range("N73").Value = -9752.467000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 1002
$ws.Range("J5").Value = 1002
$ws.Range("L5").Value = 1002
$ws.Range("N5").Value = -1226
